$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 38, column A (phone) was stored as text; the update turns it into a
# genuine numeric value.
$ws.Range("A38").Value = 71277620

# Append the new payment record (phone 71277620, Cash, 2025-08-18T17:28:56)
# as row 39, mirroring the layout of the other rows. The phone number in
# column A keeps the same (text) representation the other new-row entries
# use, so force the cell to Text before writing the numeric-looking value.
$ws.Range("A39").NumberFormat = "@"
$ws.Range("A39").Value = "71277620"
$ws.Range("C39").Value = "Cash"
$ws.Range("D39").Value = "2025-08-18T17:28:56"
$ws.Range("E39").Value = 76
$ws.Range("G39").Value = 0
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 76
